$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.174.84'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.247.37'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.02'
$ws.Range("E5").Value = '  -2.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.40'
$ws.Range("E6").Value = '  -3.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("E10").Value = '  -3.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.28'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '2.589.91'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '2.244.61'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.832'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.65'
$ws.Range("E17").Value = '  -2.63%  '
$ws.Range("D18").Value = '44.076.75'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '0.0₃0970'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.39'
$ws.Range("E20").Value = '  -3.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.37'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.54'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.42'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.95'
$ws.Range("E24").Value = '  -3.08%  '
$ws.Range("E25").Value = '  -2.35%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.99'
$ws.Range("E27").Value = '  +5.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.97'
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.11'
$ws.Range("E30").Value = '  +1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.09'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.07'
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0806'
$ws.Range("E33").Value = '  -3.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.35'
$ws.Range("E34").Value = '  +5.37%  '
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.19'
$ws.Range("E39").Value = '  -5.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.45'
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.87'
$ws.Range("E41").Value = '  -4.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0301'
$ws.Range("E42").Value = '  -3.25%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").Value = '1.735.71'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '84.98'
$ws.Range("E45").Value = '  +5.50%  '
$ws.Range("E46").Value = '  -2.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.62'
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.94'
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.10'
$ws.Range("E49").Value = '  -4.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.13'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.59'
$ws.Range("E51").Value = '  -3.88%  '
